$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 ("About Polarion QA") - Content Placeholder 2
#   - "Entirely web based" -> "Web based"
#   - "One of three main products of " run split into "One " + "of three main products of "
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange

$full3 = $tr3.Text
$oldText = "Entirely web based"
$idx = $full3.IndexOf($oldText)
if ($idx -ge 0) {
    $rng = $tr3.Characters($idx + 1, $oldText.Length)
    $rng.Text = "Web based"
}

$full3 = $tr3.Text
$splitPrefix = "One "
$idx = $full3.IndexOf("One of three main products of ")
if ($idx -ge 0) {
    $rng = $tr3.Characters($idx + 1, $splitPrefix.Length)
    $rng.Text = $splitPrefix
}

# ---------------------------------------------------------------------------
# Slide 7 ("Negatives") - Content Placeholder 2
#   - "It's costly" run split into "It's " + "costly"
#   - "Web based" paragraph replaced with "Non-tradition view"
#   - New paragraph "Manual linking" added right after it
#   - Picture repositioned
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(2)
$tr7 = $shp7.TextFrame.TextRange

$splitPrefix2 = "It" + [char]8217 + "s "
$rng = $tr7.Characters(1, $splitPrefix2.Length)
$rng.Text = $splitPrefix2

$paraWebBased = $tr7.Paragraphs(2, 1)
$paraWebBased.InsertAfter([char]13 + "Manual linking") | Out-Null

$paraWebBased2 = $tr7.Paragraphs(2, 1)
$paraWebBased2.Text = "Non-tradition view"

$pic7 = $s7.Shapes.Item(3)
$pic7.Left = 5429646 / 12700
$pic7.Top = 2865676 / 12700
